$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 49.59229566666666
$ws.Range("H2").Value2 = 148.776887
$ws.Range("I2").Value2 = 0.2250252586609286
$ws.Range("J2").Value2 = 0.2250252586609286
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.898076333333333
$ws.Range("N2").Value2 = 8.694229
$ws.Range("O2").Value2 = 0.8806895615610505
$ws.Range("P2").Value2 = 0.8806895615610505
$ws.Range("Q2").Value2 = 143.7222583872359
$ws.Range("R2").Value2 = 1293.500325485123
$ws.Range("S2").Value2 = 0.1981773963902552
$ws.Range("T2").Value2 = 0.1981773963902552

# Row 3
$ws.Range("G3").Value2 = 49.59229566666666
$ws.Range("H3").Value2 = 148.776887
$ws.Range("I3").Value2 = 0.2250252586609286
$ws.Range("J3").Value2 = 0.2250252586609286
$ws.Range("O3").Value2 = 0.05193399155394968
$ws.Range("P3").Value2 = 0.05193399155394968
$ws.Range("Q3").Value2 = 8.475257206372444
$ws.Range("R3").Value2 = 76.277314857352
$ws.Range("S3").Value2 = 0.01168645988272201
$ws.Range("T3").Value2 = 0.01168645988272201

# Row 4
$ws.Range("G4").Value2 = 49.59229566666666
$ws.Range("H4").Value2 = 148.776887
$ws.Range("I4").Value2 = 0.2250252586609286
$ws.Range("J4").Value2 = 0.2250252586609286
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.221715
$ws.Range("N4").Value2 = 0.665145
$ws.Range("O4").Value2 = 0.0673764468849998
$ws.Range("P4").Value2 = 0.0673764468849998
$ws.Range("Q4").Value2 = 10.995355833735
$ws.Range("R4").Value2 = 98.95820250361498
$ws.Range("S4").Value2 = 0.0151614023879514
$ws.Range("T4").Value2 = 0.0151614023879514

# Row 5
$ws.Range("I5").Value2 = 0.7548368219877758
$ws.Range("J5").Value2 = 0.7548368219877758
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 2.898076333333333
$ws.Range("N5").Value2 = 8.694229
$ws.Range("O5").Value2 = 0.8806895615610505
$ws.Range("P5").Value2 = 0.8806895615610505
$ws.Range("Q5").Value2 = 482.1096681123994
$ws.Range("R5").Value2 = 4338.987013011594
$ws.Range("S5").Value2 = 0.664776909806551
$ws.Range("T5").Value2 = 0.664776909806551

# Row 6
$ws.Range("I6").Value2 = 0.7548368219877758
$ws.Range("J6").Value2 = 0.7548368219877758
$ws.Range("O6").Value2 = 0.05193399155394968
$ws.Range("P6").Value2 = 0.05193399155394968
$ws.Range("S6").Value2 = 0.03920168913772337
$ws.Range("T6").Value2 = 0.03920168913772337

# Row 7
$ws.Range("I7").Value2 = 0.7548368219877758
$ws.Range("J7").Value2 = 0.7548368219877758
$ws.Range("K7").Value2 = 1
$ws.Range("L7").Value2 = 0.3333333333333333
$ws.Range("M7").Value2 = 0.221715
$ws.Range("N7").Value2 = 0.665145
$ws.Range("O7").Value2 = 0.0673764468849998
$ws.Range("P7").Value2 = 0.0673764468849998
$ws.Range("Q7").Value2 = 36.88341257133
$ws.Range("R7").Value2 = 331.95071314197
$ws.Range("S7").Value2 = 0.05085822304350143
$ws.Range("T7").Value2 = 0.05085822304350143

# Row 8
$ws.Range("G8").Value2 = 3.94583
$ws.Range("H8").Value2 = 11.83749
$ws.Range("I8").Value2 = 0.01790422089653049
$ws.Range("J8").Value2 = 0.01790422089653049
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 2.898076333333333
$ws.Range("N8").Value2 = 8.694229
$ws.Range("O8").Value2 = 0.8806895615610505
$ws.Range("P8").Value2 = 0.8806895615610505
$ws.Range("Q8").Value2 = 11.43531653835667
$ws.Range("R8").Value2 = 102.91784884521
$ws.Range("S8").Value2 = 0.01576806045145764
$ws.Range("T8").Value2 = 0.01576806045145764

# Row 9
$ws.Range("G9").Value2 = 3.94583
$ws.Range("H9").Value2 = 11.83749
$ws.Range("I9").Value2 = 0.01790422089653049
$ws.Range("J9").Value2 = 0.01790422089653049
$ws.Range("O9").Value2 = 0.05193399155394968
$ws.Range("P9").Value2 = 0.05193399155394968
$ws.Range("Q9").Value2 = 0.6743370858933334
$ws.Range("R9").Value2 = 6.069033773040001
$ws.Range("S9").Value2 = 0.0009298376568204638
$ws.Range("T9").Value2 = 0.0009298376568204638

# Row 10
$ws.Range("G10").Value2 = 3.94583
$ws.Range("H10").Value2 = 11.83749
$ws.Range("I10").Value2 = 0.01790422089653049
$ws.Range("J10").Value2 = 0.01790422089653049
$ws.Range("K10").Value2 = 1
$ws.Range("L10").Value2 = 0.3333333333333333
$ws.Range("M10").Value2 = 0.221715
$ws.Range("N10").Value2 = 0.665145
$ws.Range("O10").Value2 = 0.0673764468849998
$ws.Range("P10").Value2 = 0.0673764468849998
$ws.Range("Q10").Value2 = 0.87484969845
$ws.Range("R10").Value2 = 7.873647286050001
$ws.Range("S10").Value2 = 0.00120632278825239
$ws.Range("T10").Value2 = 0.00120632278825239

# Row 11
$ws.Range("G11").Value2 = 0.4922746666666666
$ws.Range("H11").Value2 = 1.476824
$ws.Range("I11").Value2 = 0.002233698454765135
$ws.Range("J11").Value2 = 0.002233698454765135
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 2.898076333333333
$ws.Range("N11").Value2 = 8.694229
$ws.Range("O11").Value2 = 0.8806895615610505
$ws.Range("P11").Value2 = 0.8806895615610505
$ws.Range("Q11").Value2 = 1.426649560966222
$ws.Range("R11").Value2 = 12.839846048696
$ws.Range("S11").Value2 = 0.001967194912786703
$ws.Range("T11").Value2 = 0.001967194912786703

# Row 12
$ws.Range("G12").Value2 = 0.4922746666666666
$ws.Range("H12").Value2 = 1.476824
$ws.Range("I12").Value2 = 0.002233698454765135
$ws.Range("J12").Value2 = 0.002233698454765135
$ws.Range("O12").Value2 = 0.05193399155394968
$ws.Range("P12").Value2 = 0.05193399155394968
$ws.Range("Q12").Value2 = 0.08412908416711111
$ws.Range("R12").Value2 = 0.757161757504
$ws.Range("S12").Value2 = 0.000116004876683843
$ws.Range("T12").Value2 = 0.000116004876683843

# Row 13
$ws.Range("G13").Value2 = 0.4922746666666666
$ws.Range("H13").Value2 = 1.476824
$ws.Range("I13").Value2 = 0.002233698454765135
$ws.Range("J13").Value2 = 0.002233698454765135
$ws.Range("K13").Value2 = 1
$ws.Range("L13").Value2 = 0.3333333333333333
$ws.Range("M13").Value2 = 0.221715
$ws.Range("N13").Value2 = 0.665145
$ws.Range("O13").Value2 = 0.0673764468849998
$ws.Range("P13").Value2 = 0.0673764468849998
$ws.Range("Q13").Value2 = 0.10914467772
$ws.Range("R13").Value2 = 0.9823020994799999
$ws.Range("S13").Value2 = 0.0001504986652945893
$ws.Range("T13").Value2 = 0.0001504986652945893
